$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("Croatia")

# Select the full sheet on Croatia so that once it's no longer the active tab,
# its saved view state matches Excel's default "whole sheet" selection.
[void]$src.Cells.Select()

# Duplicate the Croatia sheet (copies formatting, merged cells, column widths,
# styles, printer settings, etc.) and place it immediately after Croatia.
$src.Copy($null, $src)

# The newly inserted copy becomes the active sheet right after Croatia.
$new = $wb.Worksheets.Item($src.Index + 1)
$new.Name = "Greece"

# Update the two data cells for the new market.
$new.Range("B2").Value = "Greece Market"
$new.Range("B4").Value = "NGC-4119/T3205/T3204/T3206"

# Keep the natural "user just edited B4" selection state on the new sheet.
[void]$new.Range("B4").Select()
